# M6-RF315 BOM update: adjust panel/markup multiplier (C25) from 3 to 2.
# This ripples through every L/M column formula (L = J * $C$25, M = L * K)
# and the L25/M25 totals, removing the previously-flagged DRC quantity
# mismatches.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the driving input cell; dependent formulas auto-recalculate.
$ws.Range("C25").Value = 2

# Restore the view state captured in the saved workbook: scrolled so column
# C is the left-most visible column, with D29 as the active selection.
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$win.ScrollRow = 1
[void]$ws.Range("D29").Select()
